$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the text "60-100" in E5 with the numeric value 140
$ws.Range("E5").Value = 140

# Update the active selection from K6 to E6 (as seen in the saved workbook)
$ws.Range("E6").Select()
